# "created paint, metal livery"
#
# Appends two new documented-shader sections ("MetalColorable" and
# "MetalLiveryGloss/MetalLiveryMatte") to the end of the UV documentation,
# right after the existing "MetalChrome" section and before the section
# break (sectPr) that ends the body.
#
# We build the exact OOXML for the new paragraphs and drop it in with
# Range.InsertXML (the WordOpenXML "paste fragment" mechanism) rather than
# than driving Selection/TypeText, because that lets us reproduce the
# author's paragraph/run-property shape precisely (including the couple of
# genuinely-empty spacer paragraphs) instead of whatever defaults Word's
# typing-simulation would pick up from the caret's current formatting.

$d = $word.ActiveDocument

$newParagraphsXml = @'
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>MetalColorable</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">ShaderName: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Vehicle_Opaque_PaintGloss_Textured_LightmappedLights_Wrap</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ShaderID: BA_6C_13_00</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">TEXCOORD1: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>AoMapTextureSampler</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>LightmapLightsTextureSampler</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">TEXCOORD2: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>CrumpleTextureSampler</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ScratchTextureSampler</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>TEXCOORD6: Livery Mapping</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>MetalLiveryGloss/MetalLiveryMatte</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">ShaderName: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Vehicle_Opaque_PaintGloss_Textured_LightmappedLights_ColourOverride_Livery</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ShaderID: 56_C6_01_00</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">TEXCOORD1: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>AoMapTextureSampler</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>LightmapLightsTextureSampler</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">TEXCOORD2: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>CrumpleTextureSampler</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ScratchTextureSampler</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">TEXCOORD3: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>DiffuseTextureSampler</w:t></w:r></w:p>
'@

$packageXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + $newParagraphsXml + '</w:body>' + `
  '</w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

# Collapse to the very end of the document's main story (just before the
# final sectPr) and drop the fragment in there, after everything that is
# already in the document (the "MetalChrome" section).
$insertionPoint = $d.Content
$insertionPoint.Collapse(0)
$null = $insertionPoint.InsertXML($packageXml)
